# Update "想去人数" (interested-count) figures in the F3:F6 cells
# on both the "展览" and "全部类型" worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3" = 91
    "F4" = 49
    "F5" = 2404
    "F6" = 222
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
